$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definitions")
$ws.Activate()

# New glossary term row (appended to the "Definitions" table)
$ws.Range("A16").Value = "bruger"
$ws.Range("B16").Value = "bruger til systemet (ubestemt role) "
$ws.Range("E16").Value = "(kunden og bestillingesmodtagelse)"

# Match the look & feel of the other rows: columns A, C, D, E are centered +
# middle aligned + wrapped, column B is middle aligned + wrapped.
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("A16").VerticalAlignment = -4108
$ws.Range("A16").WrapText = $true

$ws.Range("C16").HorizontalAlignment = -4108
$ws.Range("C16").VerticalAlignment = -4108
$ws.Range("C16").WrapText = $true

$ws.Range("D16").HorizontalAlignment = -4108
$ws.Range("D16").VerticalAlignment = -4108
$ws.Range("D16").WrapText = $true

$ws.Range("E16").HorizontalAlignment = -4108
$ws.Range("E16").VerticalAlignment = -4108
$ws.Range("E16").WrapText = $true

$ws.Range("B16").VerticalAlignment = -4108
$ws.Range("B16").WrapText = $true

# Row height matches the other multi-line wrapped rows (3 lines @ 14.4pt)
$ws.Rows.Item(16).RowHeight = 43.2

# Grow the glossary table to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:E16"))

# Restore the previously-recorded selection for this sheet
$ws.Range("F14").Select() | Out-Null
